# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") sometimes holds values that look numeric (e.g. "8.00",
# "1.00", "0.0785"). Excel's Value setter auto-converts those to real numbers
# and drops the formatting (8.00 -> 8), so a leading apostrophe is used to force
# them back to literal text, exactly as the source data stores them.

# Row 2: Bitcoin
$ws.Range("D2").Value = '68.185.43'
$ws.Range("E2").Value = '  +0.58%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.638.71'
$ws.Range("E3").Value = '  +0.49%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.01%  '

# Row 5: BNB
$ws.Range("D5").Value = '''599.15'
$ws.Range("E5").Value = '  +0.51%  '

# Row 6: Solana
$ws.Range("D6").Value = '''154.60'
$ws.Range("E6").Value = '  +0.77%  '

# Row 7: USDC
$ws.Range("E7").Value = '  -0.01%  '

# Row 8: XRP
$ws.Range("D8").Value = '''0.544'
$ws.Range("E8").Value = '  -0.64%  '

# Row 9: LidoStakedEther
$ws.Range("D9").Value = '2.637.95'
$ws.Range("E9").Value = '  +0.48%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '  +7.05%  '

# Row 11: TRON
$ws.Range("E11").Value = '  -0.76%  '

# Row 12: Toncoin
$ws.Range("E12").Value = '  +0.33%  '

# Row 13: Cardano
$ws.Range("E13").Value = '  +1.47%  '

# Row 14: Avalanche
$ws.Range("D14").Value = '''27.82'
$ws.Range("E14").Value = '  +0.99%  '

# Row 15: ShibaInu
$ws.Range("E15").Value = '  +2.34%  '

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").Value = '3.116.83'
$ws.Range("E16").Value = '  +0.52%  '

# Row 17: WrappedBTC
$ws.Range("D17").Value = '68.121.25'
$ws.Range("E17").Value = '  +0.68%  '

# Row 18: WrappedEther
$ws.Range("D18").Value = '2.635.51'
$ws.Range("E18").Value = '  +0.33%  '

# Row 19: Chainlink
$ws.Range("D19").Value = '''11.35'
$ws.Range("E19").Value = '  -0.46%  '

# Row 20: BitcoinCash
$ws.Range("D20").Value = '''362.72'
$ws.Range("E20").Value = '  -1.17%  '

# Row 21: Uniswap
$ws.Range("E21").Value = '  -0.34%  '

# Row 22: Polkadot
$ws.Range("D22").Value = '''4.35'
$ws.Range("E22").Value = '  +2.89%  '

# Row 24: SuiNetwork
$ws.Range("E24").Value = '  -1.29%  '

# Row 25: Litecoin
$ws.Range("D25").Value = '''75.16'
$ws.Range("E25").Value = '  +4.46%  '

# Row 26: Dai
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  +0.01%  '

# Row 27: Aptos
$ws.Range("D27").Value = '''9.76'
$ws.Range("E27").Value = '  -0.63%  '

# Row 28: PEPE
$ws.Range("E28").Value = '  +1.70%  '

# Row 29: WrappedeETH
$ws.Range("D29").Value = '2.774.92'

# Row 30: Binance-PegBSC-USD
$ws.Range("E30").Value = '  -0.14%  '

# Row 31: Bittensor
$ws.Range("D31").Value = '''561.17'
$ws.Range("E31").Value = '  -2.53%  '

# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").Value = '''8.00'
$ws.Range("E32").Value = '  +1.24%  '

# Row 33: Fetch.AI
$ws.Range("E33").Value = '  -0.28%  '

# Row 34: PancakeSwap
$ws.Range("E34").Value = '  +0.92%  '

# Row 35: Kaspa
$ws.Range("E35").Value = '  +1.02%  '

# Row 36: FirstDigitalUSD
$ws.Range("D36").Value = '''1.00'
$ws.Range("E36").Value = '  +0.04%  '

# Row 37: ImmutableX
$ws.Range("E37").Value = '  +2.69%  '

# Row 38: Monero
$ws.Range("D38").Value = '''160.77'
$ws.Range("E38").Value = '  +0.53%  '

# Row 39: EthereumClassic
$ws.Range("E39").Value = '  +0.78%  '

# Row 41: Stacks
$ws.Range("E41").Value = '  -0.25%  '

# Row 42: RenderToken
$ws.Range("D42").Value = '''5.30'
$ws.Range("E42").Value = '  -0.67%  '

# Row 43: BabyDogeCoin
$ws.Range("D43").Value = '0.0₆0339'
$ws.Range("E43").Value = '  +1.61%  '

# Row 44: WhiteBITCoin
$ws.Range("E44").Value = '  +2.59%  '

# Row 45: dogwifhat
$ws.Range("E45").Value = '  -1.02%  '

# Row 47: OKB
$ws.Range("D47").Value = '''40.59'
$ws.Range("E47").Value = '  +1.30%  '

# Row 48: Aave
$ws.Range("D48").Value = '''157.22'
$ws.Range("E48").Value = '  +1.39%  '

# Row 50: InjectiveProtocol
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''21.75'
$ws.Range("E50").Value = '  -0.28%  '

# Row 51: Cronos
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.0785'
$ws.Range("E51").Value = '  +1.04%  '
